$d = $word.ActiveDocument

# Placeholders that get the "| nl2br " filter inserted before the closing "}}"
$names = @(
    "来店リスナー",
    "稼働",
    "推定フェーズ",
    "開催イベント",
    "良い兆候",
    "課題",
    "提案",
    "店舗様のお言葉",
    "稼働率",
    "総視聴数",
    "最大同接数",
    "UU数",
    "コミュニティいいね数",
    "ポスト合計数",
    "インプレッション数",
    "エンゲージメント数"
)

foreach ($name in $names) {
    $needle = "{{" + $name + "}}"

    $hit = $d.Content
    $found = $hit.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        continue
    }

    $matchEnd = $hit.End

    # Range covering just the trailing "}}" of the match.
    $closeRange = $d.Range($matchEnd - 2, $matchEnd)

    # Insert the filter text right before the closing braces, turning
    # "{{Name}}" into "{{Name| nl2br }}".
    $closeRange.InsertBefore("| nl2br ")

    # Toggling a character property on the newly inserted text (and
    # immediately restoring it) forces the engine to keep it as its own
    # run instead of silently re-merging it with the neighbouring runs
    # that share identical formatting - matching the 3-run split
    # ("{{Name", "| nl2br ", "}}") the target document has.
    $segStart = $matchEnd - 2
    $segEnd = $segStart + 8
    $segRange = $d.Range($segStart, $segEnd)
    $segRange.Font.Bold = 1
    $segRange.Font.Bold = 0
}
